$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# Most data rows in the income-statement table are reset to 0 (the
# "read_price" values were re-derived and zeroed out in this pass).
$zeroRows = @(11,12,13,14,16,17,19,20,21,22,24,25,26,27)
foreach ($r in $zeroRows) {
    $ws.Range("D" + $r + ":M" + $r).Value = 0
}

# Row 15: mostly "-" (not applicable), except column J which is 0.
$ws.Range("D15:I15").Value = "-"
$ws.Range("J15").Value = 0
$ws.Range("K15:M15").Value = "-"

# Row 18: columns E and F are "-", the rest are 0.
$ws.Range("D18").Value = 0
$ws.Range("E18:F18").Value = "-"
$ws.Range("G18:M18").Value = 0

# Row 23: entirely "-".
$ws.Range("D23:M23").Value = "-"
